# aggiornamento fino a 02/05
# Append daily rows 239-244 (2021-04-27 .. 2021-05-02) to the existing
# time-series table, matching the layout/formatting of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$newRows = @(
    @(239, 44313, 9,  94, 233.3780227419435),
    @(240, 44314, 3,  82, 203.5850836685039),
    @(241, 44315, 14, 78, 193.6541039773574),
    @(242, 44316, 11, 67, 166.3439098267044),
    @(243, 44317, 19, 71, 176.2748895178509),
    @(244, 44318, 14, 83, 206.0678285912905)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}

# Copy the formatting (date style on column A, plain numbers on B:D) from
# the last existing row down across the newly added ones, without
# disturbing the values just written.
$ws.Range("A238:D238").Copy() | Out-Null
$ws.Range("A239:D244").PasteSpecial(-4122) | Out-Null
